$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to stage values that Excel's smart-typing would
# otherwise misread as numbers (e.g. comma-grouped coordinate lists or
# plain decimals that must stay text, matching the source inlineStr
# cells). Writing the literal as a quoted formula keeps it text, and a
# values-only paste copies that text straight into the target cell
# without dragging a new NumberFormat/style allocation along with it.
$scratch = $ws.Cells.Item(1, 26)

function Set-TextValue($cell, $value) {
    $scratch.Formula = '="' + $value + '"'
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
    $scratch.ClearContents()
}

# Row 16
$ws.Range("D16").Value = "image_20250807111314_ppp0.jpg"
Set-TextValue $ws.Range("I16") "643,531,686,575"

# Row 17
$ws.Range("D17").Value = "image_20250807111314_ppp0.jpg"
Set-TextValue $ws.Range("I17") "794,481,830,526"
Set-TextValue $ws.Range("J17") "0.72"

# Row 18
$ws.Range("D18").Value = "image_20250808100711_ppp0.jpg"
Set-TextValue $ws.Range("I18") "1182,409,1232,451"
Set-TextValue $ws.Range("J18") "0.75"

$excel.CutCopyMode = $false
